# Applies the cryptos.xlsx price/volume refresh described by the commit
# "Updated cryptos list on Tue May 28 19:02:22 UTC 2024 with GitHub Actions".
# The Price/Volume(1h) columns hold plain-text values (the sheet stores them
# as inline strings, e.g. "1.00" / "  -3.39%  "), so before assigning each new
# value we force NumberFormat to "@" (text) to stop Excel from coercing the
# string into a Double -- which would silently drop meaningful formatting such
# as trailing zeros ("1.00" -> 1) or thousand separators ("67.577.99" would
# even fail number parsing). The style is reset back to "Normal" afterwards so
# no stray text-format style is left on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry is @(CellAddress, NewValue).
$updates = @(
    @('D2', '67.816.47'),
    @('E2', '  -3.06%  '),
    @('D3', '3.795.14'),
    @('E3', '  -3.31%  '),
    @('D4', '0.998'),
    @('E4', '  -0.36%  '),
    @('D5', '598.60'),
    @('E5', '  -1.69%  '),
    @('D6', '167.44'),
    @('E6', '  -1.62%  '),
    @('D7', '3.793.55'),
    @('E7', '  -3.43%  '),
    @('E8', '  -0.15%  '),
    @('E9', '  -2.24%  '),
    @('D10', '0.163'),
    @('E10', '  -4.07%  '),
    @('D11', '6.45'),
    @('E11', '  +0.78%  '),
    @('D12', '0.455'),
    @('E12', '  -2.96%  '),
    @('D13', '0.0000259'),
    @('E13', '  +1.11%  '),
    @('D14', '36.73'),
    @('E14', '  -4.24%  '),
    @('D15', '4.416.03'),
    @('E15', '  -3.70%  '),
    @('D16', '3.785.63'),
    @('E16', '  -3.92%  '),
    @('D17', '67.712.12'),
    @('E17', '  -3.20%  '),
    @('D18', '18.37'),
    @('E18', '  -1.77%  '),
    @('D19', '7.35'),
    @('E19', '  -3.72%  '),
    @('E20', '  -1.06%  '),
    @('D21', '10.97'),
    @('E21', '  -1.78%  '),
    @('D22', '464.45'),
    @('E22', '  -5.84%  '),
    @('D23', '0.729'),
    @('E23', '  -2.46%  '),
    @('D24', '0.0000160'),
    @('E24', '  -4.77%  '),
    @('D25', '82.45'),
    @('E25', '  -4.05%  '),
    @('D26', '2.22'),
    @('E26', '  -2.97%  '),
    @('D27', '12.02'),
    @('E27', '  -2.42%  '),
    @('E28', '  -0.20%  '),
    @('D29', '9.98'),
    @('E29', '  -1.44%  '),
    @('E30', '  -1.78%  '),
    @('D31', '3.928.93'),
    @('E31', '  -3.67%  '),
    @('D32', '7.59'),
    @('E32', '  -3.24%  '),
    @('D33', '31.25'),
    @('E33', '  -2.89%  '),
    @('D34', '2.29'),
    @('E34', '  -6.34%  '),
    @('D35', '9.44'),
    @('E35', '  -1.53%  '),
    @('D36', '3.747.52'),
    @('E36', '  -3.68%  '),
    @('D37', '0.104'),
    @('E37', '  -4.00%  '),
    @('D38', '3.62'),
    @('E38', '  +10.14%  '),
    @('D39', '0.140'),
    @('E39', '  -1.59%  '),
    @('E40', '  -3.61%  '),
    @('D41', '5.87'),
    @('E41', '  -4.41%  '),
    @('D42', '1.00'),
    @('E42', '  -0.03%  '),
    @('D43', '0.312'),
    @('E43', '  -5.40%  '),
    @('E44', '  -7.29%  '),
    @('D45', '8.69'),
    @('E45', '  +0.50%  '),
    @('B47', 'Bittensor'),
    @('C47', 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'),
    @('D47', '415.19'),
    @('E47', '  -4.87%  '),
    @('B48', 'FLOKI'),
    @('C48', 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'),
    @('D48', '141.61'),
    @('E48', '  +6.41%  '),
    @('D49', '46.80'),
    @('E49', '  -3.22%  '),
    @('D50', '141.61'),
    @('E50', '  -0.99%  '),
    @('D51', '26.01'),
    @('E51', '  +2.63%  ')
)

foreach ($update in $updates) {
    $cellRef = $update[0]
    $newValue = $update[1]
    $range = $ws.Range($cellRef)
    $range.NumberFormat = '@'
    $range.Value = $newValue
    $range.Style = 'Normal'
}
